$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.303.08"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "3.513.14"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "3.512.18"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "4.107.35"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000207"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "3.507.92"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "66.326.07"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").Value = "3.642.21"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.86%  "
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  -7.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.25"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "3.502.90"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("E37").Value = "  -3.35%  "
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0860"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.90"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("E47").Value = "  -7.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -11.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.46"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.948"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.79%  "
